$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 388
$firstRow = 2

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46074) {
        $cell.Value2 = 46075
    }
}
